$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.527.42'
$ws.Range("E2").Value = '  +0.21%  '

# Row 3
$ws.Range("D3").Value = '2.637.32'
$ws.Range("E3").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.12%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.60%  '

# Row 8
$ws.Range("E8").Value = '  +0.31%  '

# Row 9
$ws.Range("E9").Value = '  +2.50%  '

# Row 10
$ws.Range("E10").Value = '  -0.52%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.370'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.69%  '

# Row 12
$ws.Range("E12").Value = '  -0.15%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.81%  '

# Row 14
$ws.Range("D14").Value = '3.113.79'
$ws.Range("E14").Value = '  -0.24%  '

# Row 15
$ws.Range("D15").Value = '63.370.43'
$ws.Range("E15").Value = '  +0.13%  '

# Row 17
$ws.Range("D17").Value = '2.669.61'
$ws.Range("E17").Value = '  +1.14%  '

# Row 18
$ws.Range("E18").Value = '  +0.80%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.23%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.81%  '

# Row 22
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.56%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.71%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.70'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '592.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.73%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.25%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.17%  '

# Row 29
$ws.Range("E29").Value = '  -1.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.18%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.91%  '

# Row 33
$ws.Range("E33").Value = '  -2.69%  '

# Row 34
$ws.Range("E34").Value = '  +3.38%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.83%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '166.92'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.29%  '

# Row 37
$ws.Range("E37").Value = '  +1.13%  '

# Row 38
$ws.Range("E38").Value = '  -0.11%  '

# Row 39
$ws.Range("E39").Value = '  +8.67%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.12'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.41%  '

# Row 41
$ws.Range("E41").Value = '  +0.10%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '168.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.08%  '

# Row 43
$ws.Range("E43").Value = '  +1.77%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.53%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0572'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.73%  '

# Row 46
$ws.Range("E46").Value = '  +0.01%  '

# Row 47
$ws.Range("E47").Value = '  +3.82%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +15.70%  '

# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0961'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.179'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.40%  '
